$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to keep text formatting ("@") so the numeric-
# and percent-looking strings are written back as literal text, matching
# the inline-string cells already used throughout this sheet instead of
# being auto-converted to numbers/percentages by Excel.
$targetCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5",
    "D6", "E6", "D7", "D8", "E8", "D9", "E9", "D10",
    "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14",
    "E14", "D15", "E15", "E16", "D19", "E19", "E20", "D22",
    "E22", "D23", "E23", "D24", "E24", "E25", "D38", "E38",
    "D39", "E39", "D40", "E40", "E41", "D42", "E42", "D43",
    "E43", "D44", "E44", "E45", "D46", "E46", "E47", "E48",
    "E50", "E51"
)
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin price / volume figures (commit: "Updated symbol list on
# Mon Feb 13 13:48:39 UTC 2023 with GitHub Actions").
$ws.Range("D2").Value = "292.19"
$ws.Range("E2").Value = "-5.42%"
$ws.Range("D3").Value = "40.01"
$ws.Range("E3").Value = "-2.59%"
$ws.Range("D4").Value = "5.035"
$ws.Range("E4").Value = "-2.81%"
$ws.Range("D5").Value = "0.07361"
$ws.Range("E5").Value = "-3.95%"
$ws.Range("D6").Value = "4.288"
$ws.Range("E6").Value = "-0.36%"
$ws.Range("D7").Value = "1.558"
$ws.Range("D8").Value = "0.9186"
$ws.Range("E8").Value = "0.30%"
$ws.Range("D9").Value = "0.1189"
$ws.Range("E9").Value = "-4.10%"
$ws.Range("D10").Value = "0.1726"
$ws.Range("E10").Value = "-4.93%"
$ws.Range("D11").Value = "0.08733"
$ws.Range("E11").Value = "-4.55%"
$ws.Range("D12").Value = "0.04175"
$ws.Range("E12").Value = "-0.26%"
$ws.Range("D13").Value = "0.1052"
$ws.Range("E13").Value = "0.02%"
$ws.Range("D14").Value = "0.001276"
$ws.Range("E14").Value = "0.95%"
$ws.Range("D15").Value = "0.005778"
$ws.Range("E15").Value = "0.28%"
$ws.Range("E16").Value = "1.80%"
$ws.Range("D19").Value = "7.580"
$ws.Range("E19").Value = "2.68%"
$ws.Range("E20").Value = "-0.20%"
$ws.Range("D22").Value = "0.03847"
$ws.Range("E22").Value = "-4.41%"
$ws.Range("D23").Value = "0.001277"
$ws.Range("E23").Value = "0.60%"
$ws.Range("D24").Value = "0.003891"
$ws.Range("E24").Value = "-4.80%"
$ws.Range("E25").Value = "-1.78%"
$ws.Range("D38").Value = "0.02333"
$ws.Range("E38").Value = "-7.32%"
$ws.Range("D39").Value = "0.05035"
$ws.Range("E39").Value = "-4.90%"
$ws.Range("D40").Value = "0.007685"
$ws.Range("E40").Value = "-2.22%"
$ws.Range("E41").Value = "172.06%"
$ws.Range("D42").Value = "0.1269"
$ws.Range("E42").Value = "-3.04%"
$ws.Range("D43").Value = "0.007373"
$ws.Range("E43").Value = "10.59%"
$ws.Range("D44").Value = "0.007704"
$ws.Range("E44").Value = "-5.39%"
$ws.Range("E45").Value = "3.48%"
$ws.Range("D46").Value = "0.00006530"
$ws.Range("E46").Value = "-4.23%"
$ws.Range("E47").Value = "-0.22%"
$ws.Range("E48").Value = "8.63%"
$ws.Range("E50").Value = "-0.22%"
$ws.Range("E51").Value = "-0.22%"
